$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4008786678314209
$ws.Range("B1").Value = 0.3498486280441284
$ws.Range("C1").Value = 3.725913286209106
$ws.Range("D1").Value = 3.359753370285034
$ws.Range("E1").Value = 0.9117211699485779
